$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename "Active cases" -> "Active Cases"
$ws.Range("A1").Value = "Cluster name"
$ws.Range("B1").Value = "Active Cases"

# Full replacement dataset: cluster name (trailing " Outbreak" dropped) + active case
# count. Rows are alphabetically sorted by cluster name, with several clusters
# added/removed/renamed versus the previous upload.
$data = @(
    @("139 Highett St Apartment Complex Richmond", 14),
    @("3153 Sacred Heart Community St Kilda Tier 1A", 21),
    @("3528 Ottoman Village Aged Care Broadmeadows", 18),
    @("3600 Belvedere Aged Care Noble Park", 28),
    @("3612 BlueCross Glengowrie", 49),
    @("3652 Regis Aged Care Dandenong North", 21),
    @("3684 Homestyle Aged Care Langford Grange Cranbourne East", 15),
    @("3824 Estia Health South Morang", 44),
    @("3980 Arcare Keysborough Aged Care Keysborough", 25),
    @("4518 Regis Aged Care Fawkner", 21),
    @("ACFS Port Logistics Altona", 15),
    @("Alfred Health The Alfred Hospital Melbourne", 11),
    @("Armstrong Creek School Armstrong Creek", 15),
    @("Aspect Autism Spectrum Australia Disability Service Heatherton", 12),
    @("Bespoke Childcare Dingley Village", 12),
    @("Clifton Hill Primary School Clifton Hill", 13),
    @("Community Kids Pascoe Vale Early Education Centre Pascoe Vale", 23),
    @("Elements Childcare Warralily Armstrong Creek", 21),
    @("Hello Fresh Warehouse Ravenhall", 10),
    @("Honeyeater Hairdressers Bendigo", 15),
    @("Inghams Enterprises Somerville", 20),
    @("JBS Australia Brooklyn", 14),
    @("Lighthouse Christian College Cranbourne", 14),
    @("McQuinns Gym Bendigo", 22),
    @("Metcash Limited Distribution Centre Laverton North", 18),
    @("Nido Early School Wyndham Vale", 13),
    @("Northern Health The Northern Hospital Epping", 12),
    @("Pelican Childcare Craigieburn", 13),
    @("Shawlands Caravan Park Dandenong South", 15),
    @("St Vincents Hospital Emergency Department Melbourne", 24),
    @("St. Mary's Primary School Swan Hill", 11),
    @("TUROSI PTY LTD Thomastown", 13),
    @("The Robin Hood Inn Drouin West", 17),
    @("The Royal Children's Hospital Parkville", 10),
    @("Werribee Mercy Hospital Emergency Department", 40),
    @("Western Health Sunshine Hospital Emergency Department", 25),
    @("Yarrabah School Aspendale", 10),
    @("Zed Bar Albury", 13)

)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
